$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.267.30"
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("D3").Value = "1.550.31"
$ws.Range("E3").Value = "  -4.89%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.41"
$ws.Range("E5").Value = "  -3.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.00"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  -5.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0608"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.65"
$ws.Range("E10").Value = "  -4.92%  "
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("E12").Value = "  -4.88%  "
$ws.Range("D13").Value = "1.563.26"
$ws.Range("E13").Value = "  -4.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.98"
$ws.Range("E14").Value = "  -5.02%  "
$ws.Range("E15").Value = "  -4.75%  "
$ws.Range("D16").Value = "25.241.67"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").Value = "0.0₃0704"
$ws.Range("E17").Value = "  -4.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "58.46"
$ws.Range("E18").Value = "  -4.95%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "185.01"
$ws.Range("E20").Value = "  -4.12%  "
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.23"
$ws.Range("E22").Value = "  -3.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  -4.12%  "
$ws.Range("E24").Value = "  -4.38%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.82"
$ws.Range("E26").Value = "  -3.78%  "
$ws.Range("E27").Value = "  -5.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.80"
$ws.Range("E28").Value = "  -3.08%  "
$ws.Range("E29").Value = "  -5.54%  "
$ws.Range("E30").Value = "  -6.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0461"
$ws.Range("E32").Value = "  -3.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.96"
$ws.Range("E33").Value = "  -5.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.45"
$ws.Range("E34").Value = "  -3.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.32"
$ws.Range("E35").Value = "  -3.58%  "
$ws.Range("D36").Value = "1.080.37"
$ws.Range("E36").Value = "  -3.58%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0148"
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("E39").Value = "  -5.47%  "
$ws.Range("E40").Value = "  -7.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.761"
$ws.Range("E41").Value = "  -10.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.795"
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "92.49"
$ws.Range("E43").Value = "  -5.92%  "
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("D45").Value = "1.681.27"
$ws.Range("E45").Value = "  -4.84%  "
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.17"
$ws.Range("E48").Value = "  -4.30%  "
$ws.Range("E49").Value = "  -5.12%  "
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  -2.16%  "
